# Update cryptos list data (prices + 1h volume %) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.808.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.27%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.812.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.61%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'351.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.58%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'112.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +5.17%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.560"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.32%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.622"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +6.83%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'40.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +2.19%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.83%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +2.03%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'20.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.85%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.253.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.64%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.973"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +4.92%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.814.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.69%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'51.826.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +1.38%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +9.98%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.59%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'13.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.43%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.0₃0975"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.64%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'70.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.54%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'268.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.21%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +1.52%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'26.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D28").Value = "'0.162"
$ws.Range("D28").Style = "Normal"
$ws.Range("B29").Value = "'Cosmos"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'10.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +3.69%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'InjectiveProtocol"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'38.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +11.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +2.61%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'52.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.94%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'6.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.99%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.75%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0905"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +9.38%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +2.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.12%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'18.95"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +4.44%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.17%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +3.10%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.91%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.91%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Monero"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'121.71"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.81%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'WEMIXToken"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'2.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +2.23%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'EnergySwap"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'22.10"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.40%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +8.55%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +9.46%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.148.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.91%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.996"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +9.13%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +19.24%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'5.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.61%  "
$ws.Range("E51").Style = "Normal"
